$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3-12 down to 4-13.
$ws.Rows.Item(3).Insert()

# The newly inserted row 3 is blank; copy formatting/style from row 4 (old row 3) for column D (date style).
$ws.Cells.Item(4, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new weekly record in row 3.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 44764
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100114007
$ws.Cells.Item(3, 7).Value = "Jengibre"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 13000
$ws.Cells.Item(3, 13).Value = 12500
$ws.Cells.Item(3, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(3, 15).Value = "Perú"
$ws.Cells.Item(3, 16).Value = 962
$ws.Cells.Item(3, 17).Value = 13
$ws.Cells.Item(3, 18).Value = "Hortaliza"
